$wb = $excel.ActiveWorkbook

# Sheet 1: y_fitted_on_begin_2016 (col B only, rows 2-21)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 2).Value = 0.2703118165723636
$ws1.Cells.Item(3, 2).Value = 63.09324060105126
$ws1.Cells.Item(4, 2).Value = 63.85364009961007
$ws1.Cells.Item(5, 2).Value = 62.84123779559643
$ws1.Cells.Item(6, 2).Value = 61.57097913841569
$ws1.Cells.Item(7, 2).Value = 62.8117806226814
$ws1.Cells.Item(8, 2).Value = 63.28607725592121
$ws1.Cells.Item(9, 2).Value = 63.19392479836463
$ws1.Cells.Item(10, 2).Value = 62.67467080091345
$ws1.Cells.Item(11, 2).Value = 62.30377756191222
$ws1.Cells.Item(12, 2).Value = 62.99350492296069
$ws1.Cells.Item(13, 2).Value = 63.36673309872442
$ws1.Cells.Item(14, 2).Value = 63.0047660240749
$ws1.Cells.Item(15, 2).Value = 67.13935429894721
$ws1.Cells.Item(16, 2).Value = 66.00373094185461
$ws1.Cells.Item(17, 2).Value = 65.0685179682141
$ws1.Cells.Item(18, 2).Value = 65.43355294167874
$ws1.Cells.Item(19, 2).Value = 65.43775309266186
$ws1.Cells.Item(20, 2).Value = 65.17230926305706
$ws1.Cells.Item(21, 2).Value = 67.27130958490723

# Sheet 2: y_pred_on_2017_2021 (cols B, C, D, rows 2-6)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 2).Value = 68.22916511592616
$ws2.Cells.Item(2, 3).Value = 65.78527118350084
$ws2.Cells.Item(2, 4).Value = 70.67305904835149
$ws2.Cells.Item(3, 2).Value = 68.49947693249852
$ws2.Cells.Item(3, 3).Value = 65.04328898826131
$ws2.Cells.Item(3, 4).Value = 71.95566487673574
$ws2.Cells.Item(4, 2).Value = 68.76978874907088
$ws2.Cells.Item(4, 3).Value = 64.53684028980092
$ws2.Cells.Item(4, 4).Value = 73.00273720834085
$ws2.Cells.Item(5, 2).Value = 69.04010056564324
$ws2.Cells.Item(5, 3).Value = 64.15231270079259
$ws2.Cells.Item(5, 4).Value = 73.92788843049389
$ws2.Cells.Item(6, 2).Value = 69.3104123822156
$ws2.Cells.Item(6, 3).Value = 63.8456994195133
$ws2.Cells.Item(6, 4).Value = 74.77512534491791

# Sheet 3: y_fitted_on_begin_2021 (col B only, rows 2-26)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 2).Value = 0.2229949156486369
$ws3.Cells.Item(3, 2).Value = 63.04592370012753
$ws3.Cells.Item(4, 2).Value = 63.80632319868634
$ws3.Cells.Item(5, 2).Value = 62.79392089467269
$ws3.Cells.Item(6, 2).Value = 61.52366223749196
$ws3.Cells.Item(7, 2).Value = 62.76446372175767
$ws3.Cells.Item(8, 2).Value = 63.23876035499747
$ws3.Cells.Item(9, 2).Value = 63.1466078974409
$ws3.Cells.Item(10, 2).Value = 62.62735389998971
$ws3.Cells.Item(11, 2).Value = 62.25646066098849
$ws3.Cells.Item(12, 2).Value = 62.94618802203696
$ws3.Cells.Item(13, 2).Value = 63.31941619780068
$ws3.Cells.Item(14, 2).Value = 62.95744912315116
$ws3.Cells.Item(15, 2).Value = 67.09203739802349
$ws3.Cells.Item(16, 2).Value = 65.95641404093089
$ws3.Cells.Item(17, 2).Value = 65.02120106729038
$ws3.Cells.Item(18, 2).Value = 65.38623604075502
$ws3.Cells.Item(19, 2).Value = 65.39043619173815
$ws3.Cells.Item(20, 2).Value = 65.12499236213334
$ws3.Cells.Item(21, 2).Value = 67.22399268398351
$ws3.Cells.Item(22, 2).Value = 68.18184821500245
$ws3.Cells.Item(23, 2).Value = 67.33966697789542
$ws3.Cells.Item(24, 2).Value = 67.08884106839376
$ws3.Cells.Item(25, 2).Value = 67.76293986495298
$ws3.Cells.Item(26, 2).Value = 68.08037077032391

# Sheet 4: y_pred_on_2022_2026 (cols B, C, D, rows 2-6)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 2).Value = 68.39780167569482
$ws4.Cells.Item(2, 3).Value = 66.16439685499222
$ws4.Cells.Item(2, 4).Value = 70.63120649639743
$ws4.Cells.Item(3, 2).Value = 68.62079659134346
$ws4.Cells.Item(3, 3).Value = 65.4622852036364
$ws4.Cells.Item(3, 4).Value = 71.77930797905053
$ws4.Cells.Item(4, 2).Value = 68.8437915069921
$ws4.Cells.Item(4, 3).Value = 64.97542088366595
$ws4.Cells.Item(4, 4).Value = 72.71216213031826
$ws4.Cells.Item(5, 2).Value = 69.06678642264075
$ws4.Cells.Item(5, 3).Value = 64.59997678123555
$ws4.Cells.Item(5, 4).Value = 73.53359606404594
$ws4.Cells.Item(6, 2).Value = 69.28978133828939
$ws4.Cells.Item(6, 3).Value = 64.29573633792265
$ws4.Cells.Item(6, 4).Value = 74.28382633865613
